$wb = $excel.ActiveWorkbook

# --- 1) Lookup sheet: switch XMATCH() back to MATCH(...,0) per issue #3 ---
$lookup = $wb.Worksheets.Item("Lookup")
$lookup.Range("C2").Formula = "=MATCH(B2,Parameters!A3:A5,0)"
$lookup.Range("C3").Formula = "=MATCH(B3,Parameters!B2:C2,0)"
$lookup.Range("C2").Select() | Out-Null

# --- 2) Model sheet: change Storage type selection from Tank to Lagoon ---
$model = $wb.Worksheets.Item("Model")
$model.Range("C14").Value = "Lagoon"

# --- 3) ChangeLog sheet: add new version row documenting the change ---
$cl = $wb.Worksheets.Item("ChangeLog")
$cl.Range("A11").Value = 1.3
$cl.Range("B11").NumberFormat = $cl.Range("B10").NumberFormat
$cl.Range("B11").Value = 45156
$cl.Range("C11").Value = "AMOSTO.xlsx"
$cl.Range("D11").Value = "Sasha"
$cl.Range("E11").Value = "Switch back to MATCH(" + [char]0x2026 + ",0). See issue #3."
$cl.Range("E12").Select() | Out-Null
